# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect the latest scrape numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 875
    6  = 341
    7  = 10648
    8  = 250
    10 = 4
    12 = 138
    16 = 39
    20 = 1033
    22 = 103
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
